$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.624.03'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.596.40'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.12'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.51'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '1.819.87'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '1.606.84'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.07'
$ws.Range("E16").Value = '  +0.05%  '
$ws.Range("D17").Value = '26.604.02'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '0.0₃0738'
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.45'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.03'
$ws.Range("E21").Value = '  +4.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.89'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.12'
$ws.Range("E25").Value = '  -1.06%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0509'
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.94'
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").Value = '1.278.62'
$ws.Range("E34").Value = '  -1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.619'
$ws.Range("E35").Value = '  -7.70%  '
$ws.Range("E36").Value = '  +0.87%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.839'
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  +17.98%  '
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.18'
$ws.Range("E43").Value = '  +0.63%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.784'
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").Value = '1.732.69'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.18'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("E48").Value = '  +3.88%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0508'
$ws.Range("E49").Value = '  +0.88%  '
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.44'
$ws.Range("E51").Value = '  -0.70%  '
